$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (historical_growth_revenue_last_5_years) is no longer populated for rows 2-3
$ws.Range("D2:D3").ClearContents()

# Row 2: refresh capital-structure / margin metrics
$ws.Range("G2").Value = -0.2850617283950617
$ws.Range("H2").Value = -0.3037037037037037
$ws.Range("I2").Value = -0.2695473251028807
$ws.Range("J2").Value = -0.2695473251028807
$ws.Range("K2").Value = -6.55
$ws.Range("L2").Value = -0.2695473251028807
$ws.Range("U2").Value = 0.818
$ws.Range("V2").Value = 0.01455516014234875
$ws.Range("W2").Value = -0.8675496688741722
$ws.Range("X2").Value = 0.06424379446561708
$ws.Range("Y2").Value = -0.9317934633397893
$ws.Range("Z2").Value = 4.942037827943867
$ws.Range("AA2").Value = -1.33211307707952
$ws.Range("AB2").Value = 0.06312149333162342
$ws.Range("AC2").Value = -1.395234570411143
$ws.Range("AD2").Value = 1.69
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1.69
$ws.Range("AG2").Value = 0.872
$ws.Range("AH2").Value = 0.02919329763344274
$ws.Range("AI2").Value = 0.4870317002881844
$ws.Range("AJ2").Value = 0.0152789458929072
$ws.Range("AK2").Value = 0.3288084464555053
$ws.Range("AL2").Value = 0.22
$ws.Range("AM2").Value = 0.22
$ws.Range("AN2").Value = -0.2522388059701492
$ws.Range("AO2").Value = -29.77272727272727
$ws.Range("AP2").Value = -0.1301492537313433
$ws.Range("AQ2").Value = -29.77272727272727

# Row 3: refresh capital-structure / margin metrics
$ws.Range("G3").Value = -0.2850617283950617
$ws.Range("H3").Value = -0.3037037037037037
$ws.Range("I3").Value = -0.2695473251028807
$ws.Range("J3").Value = -0.2695473251028807
$ws.Range("K3").Value = -6.55
$ws.Range("L3").Value = -0.2695473251028807
$ws.Range("U3").Value = 0.818
$ws.Range("V3").Value = 0.01455516014234875
$ws.Range("W3").Value = -0.8675496688741722
$ws.Range("X3").Value = 0.06424379446561708
$ws.Range("Y3").Value = -0.9317934633397893
$ws.Range("Z3").Value = 4.942037827943867
$ws.Range("AA3").Value = -1.33211307707952
$ws.Range("AB3").Value = 0.06312149333162342
$ws.Range("AC3").Value = -1.395234570411143
$ws.Range("AD3").Value = 1.69
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.69
$ws.Range("AG3").Value = 0.872
$ws.Range("AH3").Value = 0.02919329763344274
$ws.Range("AI3").Value = 0.4870317002881844
$ws.Range("AJ3").Value = 0.0152789458929072
$ws.Range("AK3").Value = 0.3288084464555053
$ws.Range("AL3").Value = 0.22
$ws.Range("AM3").Value = 0.22
$ws.Range("AN3").Value = -0.2522388059701492
$ws.Range("AO3").Value = -29.77272727272727
$ws.Range("AP3").Value = -0.1301492537313433
$ws.Range("AQ3").Value = -29.77272727272727

